# fcc-cube_AuCu_data.xlsx - update CE (column F) values, night of 02.25.2019 re-run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F23").Value = -3.074837574169006
$ws.Range("F24").Value = -3.086218696333423
$ws.Range("F26").Value = -3.107829617910939
$ws.Range("F27").Value = -3.118059417324035
$ws.Range("F28").Value = -3.1236713448568
$ws.Range("F29").Value = -3.131844671775884
$ws.Range("F30").Value = -3.139513071802657
$ws.Range("F31").Value = -3.145124999335423
$ws.Range("F32").Value = -3.148680454374174
$ws.Range("F33").Value = -3.156348854400949
$ws.Range("F34").Value = -3.160550705298712
$ws.Range("F35").Value = -3.167572709466471
$ws.Range("F36").Value = -3.171774560364238
$ws.Range("F37").Value = -3.175976411262003
$ws.Range("F38").Value = -3.180178262159763
$ws.Range("F39").Value = -3.18579018969253
$ws.Range("F40").Value = -3.188581963955295
$ws.Range("F41").Value = -3.194193891488059
$ws.Range("F42").Value = -3.198395742385824
$ws.Range("F43").Value = -3.19977744001359
$ws.Range("F45").Value = -3.208181141809117
$ws.Range("F46").Value = -3.180166354389811
$ws.Range("F47").Value = -3.199231787493216
$ws.Range("F48").Value = -3.184339826280342
$ws.Range("F49").Value = -3.220758315495175
$ws.Range("F50").Value = -3.189383676170477
$ws.Range("F52").Value = -3.179317770846336
$ws.Range("F53").Value = -3.164002907064864
$ws.Range("F54").Value = -3.168087961895831
$ws.Range("F55").Value = -3.152643721172519
$ws.Range("F56").Value = -3.152627923040326
$ws.Range("F57").Value = -3.151975339042123
$ws.Range("F58").Value = -3.139988037454823
$ws.Range("F59").Value = -3.124556377606552
$ws.Range("F60").Value = -3.105222655219386
$ws.Range("F61").Value = -3.08588893283222
$ws.Range("F62").Value = -3.066555210445059
$ws.Range("F64").Value = -3.027887765670724
$ws.Range("F71").Value = -2.892551708960567
$ws.Range("F81").Value = -3.308740685661148
$ws.Range("F82").Value = -3.361829481115791
$ws.Range("F83").Value = -3.38462851348031
$ws.Range("F84").Value = -3.361525338348633
$ws.Range("F85").Value = -3.329604919310752
$ws.Range("F86").Value = -3.294181008040878
$ws.Range("F87").Value = -3.24245275237227
$ws.Range("F88").Value = -3.16934638709593
$ws.Range("F92").Value = -3.407583622036351
$ws.Range("F93").Value = -3.441862792427939
$ws.Range("F94").Value = -3.449057862828091
$ws.Range("F95").Value = -3.463446714426212
$ws.Range("F96").Value = -3.438447729642274
$ws.Range("F97").Value = -3.408477220061477
$ws.Range("F98").Value = -3.345664232817219
$ws.Range("F99").Value = -3.278589661508707
$ws.Range("F100").Value = -3.191750915264368
$ws.Range("F104").Value = -3.507385378780688
$ws.Range("F105").Value = -3.512671387526732
$ws.Range("F106").Value = -3.511617738984432
$ws.Range("F107").Value = -3.492194608629767
$ws.Range("F108").Value = -3.454151431462957
$ws.Range("F109").Value = -3.403003174993025
$ws.Range("F110").Value = -3.337792523283967
$ws.Range("F111").Value = -3.254468085314206
$ws.Range("F114").Value = -3.530024423105069
$ws.Range("F115").Value = -3.551351387092928
$ws.Range("F116").Value = -3.559360054499963
$ws.Range("F117").Value = -3.553275620643342
$ws.Range("F118").Value = -3.533364998764425
$ws.Range("F119").Value = -3.500632984990073
$ws.Range("F120").Value = -3.447201709007445
$ws.Range("F121").Value = -3.381872592191481
$ws.Range("F122").Value = -3.304214849129697
$ws.Range("F125").Value = -3.568038739141522
$ws.Range("F126").Value = -3.587744679768577
$ws.Range("F127").Value = -3.595768758865329
$ws.Range("F128").Value = -3.589219997587994
$ws.Range("F129").Value = -3.563351863010789
$ws.Range("F130").Value = -3.526873198727701
$ws.Range("F131").Value = -3.478116458038564
$ws.Range("F132").Value = -3.412872959626577
$ws.Range("F133").Value = -3.335024735037095
$ws.Range("F136").Value = -3.596654125559561
$ws.Range("F137").Value = -3.611794805617069
$ws.Range("F138").Value = -3.613171719522013
$ws.Range("F139").Value = -3.612669629626544
$ws.Range("F140").Value = -3.582301255950159
$ws.Range("F141").Value = -3.554852134860904
$ws.Range("F142").Value = -3.504367203980358
$ws.Range("F143").Value = -3.439801387576286
$ws.Range("F144").Value = -3.367652146589238
$ws.Range("F147").Value = -3.620088982644755
$ws.Range("F148").Value = -3.636444531301638
$ws.Range("F149").Value = -3.640793682156705
$ws.Range("F150").Value = -3.630597738556407
$ws.Range("F151").Value = -3.606883191810213
$ws.Range("F152").Value = -3.569309448101849
$ws.Range("F153").Value = -3.523687848075992
$ws.Range("F154").Value = -3.457569927260804
$ws.Range("F155").Value = -3.386610986009253
$ws.Range("F158").Value = -3.638867878011204
$ws.Range("F159").Value = -3.650310193496958
$ws.Range("F160").Value = -3.655165881137822
$ws.Range("F161").Value = -3.644795729929751
$ws.Range("F162").Value = -3.622423535095442
$ws.Range("F163").Value = -3.58561681049506
$ws.Range("F164").Value = -3.530650288464662
$ws.Range("F165").Value = -3.472505963866745
$ws.Range("F166").Value = -3.404109092308625
$ws.Range("F169").Value = -3.65316996447408
$ws.Range("F170").Value = -3.667953355227468
$ws.Range("F171").Value = -3.669480757653198
$ws.Range("F172").Value = -3.662842411099027
$ws.Range("F173").Value = -3.635261251748638
$ws.Range("F174").Value = -3.596856620482297
$ws.Range("F175").Value = -3.545644036104455
$ws.Range("F176").Value = -3.485066444501304
$ws.Range("F177").Value = -3.417819696198888
$ws.Range("F180").Value = -3.666708727290143
$ws.Range("F181").Value = -3.677601946326233
$ws.Range("F182").Value = -3.677732539760361
$ws.Range("F183").Value = -3.667091882604479
$ws.Range("F184").Value = -3.64413181541423
$ws.Range("F185").Value = -3.607821206513041
$ws.Range("F186").Value = -3.555663589615599
$ws.Range("F187").Value = -3.494751408404368
$ws.Range("F188").Value = -3.429630506479737
$ws.Range("F191").Value = -3.677621181168572
$ws.Range("F192").Value = -3.689716603712566
$ws.Range("F193").Value = -3.69172456024225
$ws.Range("F194").Value = -3.679214905547041
$ws.Range("F195").Value = -3.652983292794469
$ws.Range("F196").Value = -3.616856413706743
$ws.Range("F197").Value = -3.566305265202362
$ws.Range("F198").Value = -3.50393204074598
$ws.Range("F199").Value = -3.439384557014623
$ws.Range("F202").Value = -3.687252570886743
$ws.Range("F203").Value = -3.69824074727198
$ws.Range("F204").Value = -3.697848009955548
$ws.Range("F205").Value = -3.684313369368108
$ws.Range("F206").Value = -3.659103725647463
$ws.Range("F207").Value = -3.62383874320135
$ws.Range("F208").Value = -3.574393347550454
$ws.Range("F209").Value = -3.511169419601352
$ws.Range("F210").Value = -3.447776083538979
